$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.252.23"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "3.424.38"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "413.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.627"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.15%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -2.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.140"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.76"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000220"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.85%  "
$ws.Range("D14").Value = "3.962.47"
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.79%  "
$ws.Range("D17").Value = "3.425.38"
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.70"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.70%  "
$ws.Range("E19").Value = "  -1.58%  "
$ws.Range("D20").Value = "62.244.37"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "473.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.61%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("E23").Value = "  +2.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +10.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "33.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.75%  "
$ws.Range("E28").Value = "  +0.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.54%  "
$ws.Range("E30").Value = "  -1.50%  "
$ws.Range("E31").Value = "  -3.08%  "
$ws.Range("E32").Value = "  -1.74%  "
$ws.Range("E33").Value = "  -3.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "40.88"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.01%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.97"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0487"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.03"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.74%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("E41").Value = "  +1.88%  "
$ws.Range("B42").Value = "LidoDAOToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.32"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.14%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "145.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.65"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.33"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.06%  "
$ws.Range("E46").Value = "  +3.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.36"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +17.62%  "
$ws.Range("E48").Value = "  -2.55%  "
$ws.Range("D49").Value = "0.0₃0546"
$ws.Range("E49").Value = "  +26.40%  "
$ws.Range("E50").Value = "  -0.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "112.78"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.07%  "
